$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RACP")
$ws.Range("B2").Formula = "=ROUND(100/About!A11,0)"
